{"js": "// Replace the two-digit multiplication problems in the worksheet table\n// with their updated values, as described by the commit diff.\n// Each \"old\" expression is unique within the document, so a straight\n// search-and-replace of each pair is safe and order-independent.\nconst replacements = [\n  [\"60\u00d784=\", \"56\u00d731=\"],\n  [\"66\u00d731=\", \"70\u00d711=\"],\n  [\"66\u00d761=\", \"73\u00d732=\"],\n  [\"67\u00d759=\", \"89\u00d726=\"],\n  [\"69\u00d787=\", \"83\u00d731=\"],\n  [\"72\u00d712=\", \"20\u00d711=\"],\n  [\"34\u00d744=\", \"75\u00d749=\"],\n  [\"24\u00d780=\", \"16\u00d713=\"],\n  [\"92\u00d792=\", \"79\u00d767=\"],\n  [\"53\u00d749=\", \"56\u00d772=\"],\n  [\"44\u00d715=\", \"66\u00d729=\"],\n  [\"65\u00d754=\", \"16\u00d711=\"],\n  [\"13\u00d718=\", \"17\u00d762=\"],\n  [\"89\u00d757=\", \"88\u00d713=\"],\n  [\"28\u00d744=\", \"19\u00d734=\"],\n  [\"26\u00d761=\", \"26\u00d757=\"],\n  [\"15\u00d799=\", \"54\u00d754=\"],\n  [\"61\u00d755=\", \"30\u00d716=\"],\n  [\"40\u00d751=\", \"47\u00d733=\"],\n  [\"60\u00d763=\", \"16\u00d796=\"],\n  [\"61\u00d783=\", \"39\u00d788=\"],\n  [\"97\u00d712=\", \"13\u00d771=\"],\n  [\"42\u00d734=\", \"56\u00d749=\"],\n  [\"14\u00d759=\", \"16\u00d739=\"],\n  [\"52\u00d733=\", \"48\u00d758=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the worksheet table\n# with their updated values, as described by the commit diff.\n# Each \"old\" expression is unique within the document, so a straight\n# find-and-replace of each pair (scoped to the whole document body) is\n# safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"60\u00d784=\", \"56\u00d731=\"),\n    @(\"66\u00d731=\", \"70\u00d711=\"),\n    @(\"66\u00d761=\", \"73\u00d732=\"),\n    @(\"67\u00d759=\", \"89\u00d726=\"),\n    @(\"69\u00d787=\", \"83\u00d731=\"),\n    @(\"72\u00d712=\", \"20\u00d711=\"),\n    @(\"34\u00d744=\", \"75\u00d749=\"),\n    @(\"24\u00d780=\", \"16\u00d713=\"),\n    @(\"92\u00d792=\", \"79\u00d767=\"),\n    @(\"53\u00d749=\", \"56\u00d772=\"),\n    @(\"44\u00d715=\", \"66\u00d729=\"),\n    @(\"65\u00d754=\", \"16\u00d711=\"),\n    @(\"13\u00d718=\", \"17\u00d762=\"),\n    @(\"89\u00d757=\", \"88\u00d713=\"),\n    @(\"28\u00d744=\", \"19\u00d734=\"),\n    @(\"26\u00d761=\", \"26\u00d757=\"),\n    @(\"15\u00d799=\", \"54\u00d754=\"),\n    @(\"61\u00d755=\", \"30\u00d716=\"),\n    @(\"40\u00d751=\", \"47\u00d733=\"),\n    @(\"60\u00d763=\", \"16\u00d796=\"),\n    @(\"61\u00d783=\", \"39\u00d788=\"),\n    @(\"97\u00d712=\", \"13\u00d771=\"),\n    @(\"42\u00d734=\", \"56\u00d749=\"),\n    @(\"14\u00d759=\", \"16\u00d739=\"),\n    @(\"52\u00d733=\", \"48\u00d758=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
